# Apply weekly billing update to the "WEEKLY UNITS COMPLETED PER SCOPE ID" report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: regenerate timestamp ---
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# --- Report summary: total billed amount now populated ---
$ws.Range("C8").Value = 7309.41

# --- Report details: clear out the Scope ID value ---
$ws.Range("G10").Value = ""

# --- Per line-item pricing (column H) for the detail rows ---
$ws.Range("H16").Value = 648.53
$ws.Range("H17").Value = 648.53
$ws.Range("H18").Value = 648.53
$ws.Range("H19").Value = 648.53
$ws.Range("H20").Value = 648.53
$ws.Range("H21").Value = 648.53
$ws.Range("H22").Value = 198.88
$ws.Range("H23").Value = 478.55
$ws.Range("H24").Value = 198.88
$ws.Range("H25").Value = 478.55
$ws.Range("H26").Value = 198.88
$ws.Range("H27").Value = 478.55
$ws.Range("H28").Value = 198.88
$ws.Range("H29").Value = 478.55
$ws.Range("H30").Value = 198.88
$ws.Range("H31").Value = 478.55
$ws.Range("H32").Value = 31.08

# --- TOTAL row ---
$ws.Range("H33").Value = 7309.410000000001
